# Add a new worksheet "ArchivedSamples" at the end of the workbook, populate
# it with the same EndPoint/gridName/filters layout used by the other
# "grid" sheets, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Create the new sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ArchivedSamples"

# Column widths similar to the other sheets.
$newSheet.Columns.Item(1).ColumnWidth = 24
$newSheet.Columns.Item(2).ColumnWidth = 25.7109375

# Header row (merged A1:B1) -- "Assert200"
$newSheet.Range("A1:B1").Merge()
$newSheet.Range("A1").Value = "Assert200"

# Row 2 -- column titles
$newSheet.Range("A2").Value = "EndPoint"
$newSheet.Range("B2").Value = "gridName"

# Row 3 -- endpoint values
$newSheet.Range("A3").Value = "/gridLayout/filters"
$newSheet.Range("B3").Value = "ArchiveSamples"

# Row 4 left empty on purpose.

# Selection/active cell for the new sheet.
$newSheet.Range("J31").Select()

# Make the new sheet the active tab.
$newSheet.Activate()

# Update the selection on the DonationGrid sheet (second sheet).
$donationSheet = $wb.Worksheets.Item(2)
$donationSheet.Range("A1:B4").Select()
